# Apply the workbook edit described in the diff:
# 1) Populate the previously-empty sheet "15_" with a new RC-circuit question.
# 2) Add a new sheet "16_" at the end of the workbook with a milk/capacitor comparison question.
# 3) Update the active sheet / selections to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "15_" - fill in the RC circuit / capacitor voltage question
# ---------------------------------------------------------------------
$s15 = $wb.Worksheets.Item("15_")

$s15.Range("A1").Value = 'What causes the voltage at the top capacitor plate to increase when the battery is attached to the RC circuit shown above?  Mark all that are true'
$s15.Rows.Item(1).RowHeight = 75

$s15.Range("A2").Value = 'The high voltage at the positive battery terminal drives electrons through the resistor towards  the capacitor plate'
$s15.Range("B2").Value = 'Y'
$s15.Range("C2").Value = 'Yep!  The positive battery voltage will always be higher or equal to the top capacitor plate voltage, so electrons will only move from left to right through the resistor'
$s15.Rows.Item(2).RowHeight = 60

$s15.Range("A3").Value = 'The resistor pushes electrons from left to right towards the capacitor plate'
$s15.Range("B3").Value = 'N'
$s15.Rows.Item(3).RowHeight = 45

$s15.Range("A4").Value = 'The top capacitor plate stays at zero voltage which attracts electrons'
$s15.Range("B4").Value = 'N'
$s15.Rows.Item(4).RowHeight = 45

$s15.Range("A5").Value = 'Electrons are attracted to the zero voltage of the lower plate but can''t cross the gap because there is no conductive wire'
$s15.Range("B5").Value = 'Y'
$s15.Range("C5").Value = 'Yep!  The electrons want to go to "ground" but can''t because of the gap between the plates.'
$s15.Rows.Item(5).RowHeight = 60

# Extend formatting (wrap text) and approximate column widths down through row 15
# to mirror the template used by the other question sheets in this workbook.
$s15.Range("A1:C15").WrapText = $true
$s15.Columns.Item(1).ColumnWidth = 29.03
$s15.Columns.Item(3).ColumnWidth = 42.6

$s15.Range("C3").Select()

# ---------------------------------------------------------------------
# 2) Add new sheet "16_" at the end of the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$s16 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$s16.Name = "16_"

$s16.Range("A1").Value = 'The plot of the voltage at the capacitor plate looks like the change in the temperature of the milk returning to the environmental temperature.  How might the two things be similar?  (Mark all that are true)'
$s16.Rows.Item(1).RowHeight = 105

$s16.Range("A2").Value = 'The two situations are governed by a similar diffrential equation'
$s16.Range("B2").Value = 'Y'
$s16.Range("C2").Value = 'Yep!  Surprising, huh?'
$s16.Rows.Item(2).RowHeight = 30

$s16.Range("A3").Value = 'The rate of change in each situation is determined by how far the state variable is from some steady state situation.'
$s16.Range("B3").Value = 'Y'
$s16.Range("C3").Value = 'Yep!  This is the heart of the similarity: the milk heats more rapidly when its temperature is far from the environmental temperature, while the capacitor fills with electrons more rapidly when its voltage is much different than the battery voltage.'
$s16.Rows.Item(3).RowHeight = 135

$s16.Range("A4").Value = 'Milk is actually one kind of capacitor'
$s16.Range("B4").Value = 'N'
$s16.Range("C4").Value = 'Sadly, no.  These are two totally different physical systems that behave mathematically in a very similar way.'
$s16.Rows.Item(4).RowHeight = 75

$s16.Range("A5").Value = 'Both situations could be described as experiencing an "exponential decline"'
$s16.Range("B5").Value = 'Y'
$s16.Range("C5").Value = 'Yep!  Even though both states are going "up", they both move exponentially toward some steady value, and so we can think of this as a "decline"'
$s16.Rows.Item(5).RowHeight = 75

# Extend formatting (wrap text) and approximate column widths down through row 16
$s16.Range("A1:E16").WrapText = $true
$s16.Columns.Item(1).ColumnWidth = 29.6
$s16.Columns.Item(3).ColumnWidth = 27.88

# Make "16_" the active sheet/tab, with the diff-specified selection
$s16.Activate()
$s16.Range("E3").Select()
